# Update "想去人数" (F column) figures on the sheets that carry the
# full event list ("展览" and "全部类型"). Both sheets previously had
# identical F-column values and both receive the same updates.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 7541
    5  = 322
    6  = 29
    7  = 18
    8  = 22
    9  = 5654
    10 = 143
    12 = 17
    13 = 1739
    15 = 1212
    16 = 275
    17 = 5495
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
